$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("101_2")

# Correct the "End Date" row (row 5) values which were mistakenly entered
# as 10/28/1993 (serial 34270). They should be 10/28/1990 (serial 33174).
$ws.Range("B5:D5").Value = 33174
